{"js": "// The edit removes the explicit figure/table numbers from two in-text\n// references (\"kuvion 1 mukaisesti:\" -> \"kuvion mukaisesti:\" and\n// \"...taulukossa 1:\" -> \"...taulukossa :\"), while leaving the actual\n// figure/table captions (\"Kuvio 1: ...\", \"Taulukko 1: ...\") untouched.\n\n// 1) \"kuvion 1 mukaisesti:\" -> \"kuvion mukaisesti:\"\nconst figRefs = context.document.body.search(\"kuvion 1 mukaisesti:\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nfigRefs.load(\"text\");\nawait context.sync();\n\nfor (const r of figRefs.items) {\n  r.insertText(\"kuvion mukaisesti:\", \"Replace\");\n}\nawait context.sync();\n\n// 2) \"m\u00e4\u00e4r\u00e4 on esitetty taulukossa 1:\" -> \"m\u00e4\u00e4r\u00e4 on esitetty taulukossa :\"\nconst tblRefs = context.document.body.search(\n  \"m\u00e4\u00e4r\u00e4 on esitetty taulukossa 1:\",\n  { matchCase: true, matchWholeWord: false }\n);\ntblRefs.load(\"text\");\nawait context.sync();\n\nfor (const r of tblRefs.items) {\n  r.insertText(\"m\u00e4\u00e4r\u00e4 on esitetty taulukossa :\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The edit removes the explicit figure/table numbers from two in-text\n# references:\n#   \"kuvion 1 mukaisesti:\"              -> \"kuvion mukaisesti:\"\n#   \"...m\u00e4\u00e4r\u00e4 on esitetty taulukossa 1:\" -> \"...m\u00e4\u00e4r\u00e4 on esitetty taulukossa :\"\n# The actual figure/table captions (\"Kuvio 1: ...\", \"Taulukko 1: ...\")\n# are left untouched.\n\n$d = $word.ActiveDocument\n\n# 1) \"kuvion 1 mukaisesti:\" -> \"kuvion mukaisesti:\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"kuvion 1 mukaisesti:\"\n$find1.Replacement.Text = \"kuvion mukaisesti:\"\n$find1.Forward = $true\n$find1.Wrap = 1\n$find1.Format = $false\n$find1.MatchCase = $true\n$find1.MatchWholeWord = $false\n$find1.MatchWildcards = $false\n$find1.Execute(\n    [ref]$find1.Text,\n    [ref]$find1.MatchCase,\n    [ref]$find1.MatchWholeWord,\n    [ref]$find1.MatchWildcards,\n    $null,\n    $null,\n    [ref]$find1.Forward,\n    [ref]$find1.Wrap,\n    $null,\n    [ref]$find1.Replacement.Text,\n    2\n)\n\n# 2) \"m\u00e4\u00e4r\u00e4 on esitetty taulukossa 1:\" -> \"m\u00e4\u00e4r\u00e4 on esitetty taulukossa :\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"m\u00e4\u00e4r\u00e4 on esitetty taulukossa 1:\"\n$find2.Replacement.Text = \"m\u00e4\u00e4r\u00e4 on esitetty taulukossa :\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Format = $false\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.MatchWildcards = $false\n$find2.Execute(\n    [ref]$find2.Text,\n    [ref]$find2.MatchCase,\n    [ref]$find2.MatchWholeWord,\n    [ref]$find2.MatchWildcards,\n    $null,\n    $null,\n    [ref]$find2.Forward,\n    [ref]$find2.Wrap,\n    $null,\n    [ref]$find2.Replacement.Text,\n    2\n)\n"}
